# Refresh the cryptocurrency price / 1h-volume table with newly scraped
# values (GitHub Actions cron update). Price cells that look like plain
# decimals are pre-formatted as Text ("@") so Excel keeps them as strings
# (matching the sheet's existing inline-string convention) instead of
# auto-converting them to numbers; prices that already contain extra dots
# (e.g. "64.357.37") are unambiguous and don't need that treatment.
# Rows 37/38 also swap places: dogwifhat now outranks Bittensor.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.357.37"
$ws.Range("E2").Value = "  +1.80%  "
$ws.Range("D3").Value = "3.084.72"
$ws.Range("E3").Value = "  +0.93%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "560.82"
$ws.Range("E5").Value = "  +2.20%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.75"
$ws.Range("E6").Value = "  +4.31%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").Value = "3.082.60"
$ws.Range("E8").Value = "  +1.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.507"
$ws.Range("E9").Value = "  +1.19%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.154"
$ws.Range("E10").Value = "  +2.40%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.16"
$ws.Range("E11").Value = "  -3.61%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.472"
$ws.Range("E12").Value = "  +4.43%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000229"
$ws.Range("E13").Value = "  +1.25%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.28"
$ws.Range("E14").Value = "  +1.72%  "
$ws.Range("D15").Value = "3.578.76"
$ws.Range("E15").Value = "  +0.87%  "
$ws.Range("D16").Value = "64.344.87"
$ws.Range("E16").Value = "  +1.69%  "
$ws.Range("D17").Value = "3.080.34"
$ws.Range("E17").Value = "  +0.78%  "
$ws.Range("E18").Value = "  +1.35%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.78"
$ws.Range("E19").Value = "  +0.80%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "480.69"
$ws.Range("E20").Value = "  -0.03%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.94"
$ws.Range("E21").Value = "  +2.09%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.677"
$ws.Range("E22").Value = "  +0.86%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.57"
$ws.Range("E23").Value = "  +5.22%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.85"
$ws.Range("E24").Value = "  +10.92%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "81.30"
$ws.Range("E25").Value = "  +0.90%  "
$ws.Range("E26").Value = "  -0.07%  "
$ws.Range("E27").Value = "  +2.52%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.06"
$ws.Range("E28").Value = "  +2.04%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.09"
$ws.Range("E29").Value = "  +5.47%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.998"
$ws.Range("E30").Value = "  -0.22%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "26.25"
$ws.Range("E31").Value = "  +1.18%  "
$ws.Range("E32").Value = "  -0.17%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.51"
$ws.Range("E33").Value = "  +3.51%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.60"
$ws.Range("E34").Value = "  -1.63%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.21"
$ws.Range("E35").Value = "  +4.16%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "55.80"
$ws.Range("E36").Value = "  +0.75%  "
$ws.Range("B37").Value = "dogwifhat"
$ws.Range("C37").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.05"
$ws.Range("E37").Value = "  +18.34%  "
$ws.Range("B38").Value = "Bittensor"
$ws.Range("C38").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "458.51"
$ws.Range("E38").Value = "  -1.31%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0408"
$ws.Range("E39").Value = "  +3.24%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0827"
$ws.Range("E40").Value = "  +1.59%  "
$ws.Range("D41").Value = "2.977.21"
$ws.Range("E41").Value = "  -2.91%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.27"
$ws.Range("E42").Value = "  +0.45%  "
$ws.Range("E43").Value = "  -1.80%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "27.98"
$ws.Range("E44").Value = "  -1.05%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.262"
$ws.Range("E45").Value = "  +3.84%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.16"
$ws.Range("E46").Value = "  +5.78%  "
$ws.Range("E47").Value = "  +0.02%  "
$ws.Range("E48").Value = "  +2.68%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "120.99"
$ws.Range("E49").Value = "  +3.65%  "
$ws.Range("D50").Value = "0.0₃0518"
$ws.Range("E50").Value = "  +1.84%  "
$ws.Range("E51").Value = "  +1.08%  "
